# Add a new worksheet "white tags" (with RGB sample data + min/max formulas)
# after the existing two sheets, mirroring the target commit
# "added sheet for rgb values of the white tags".

$wb = $excel.ActiveWorkbook

# New sheet goes after the last existing sheet (so it lands in position 3).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "white tags"

# Header row
$ws.Range("A1").Value = "R"
$ws.Range("B1").Value = "G"
$ws.Range("C1").Value = "B"

# RGB sample rows (2-34)
$ws.Range("A2").Value = 255
$ws.Range("B2").Value = 255
$ws.Range("C2").Value = 253

$ws.Range("A3").Value = 247
$ws.Range("B3").Value = 247
$ws.Range("C3").Value = 247

$ws.Range("A4").Value = 252
$ws.Range("B4").Value = 254
$ws.Range("C4").Value = 251

$ws.Range("A5").Value = 252
$ws.Range("B5").Value = 251
$ws.Range("C5").Value = 247

$ws.Range("A6").Value = 250
$ws.Range("B6").Value = 248
$ws.Range("C6").Value = 235

$ws.Range("A7").Value = 254
$ws.Range("B7").Value = 255
$ws.Range("C7").Value = 255

$ws.Range("A8").Value = 233
$ws.Range("B8").Value = 207
$ws.Range("C8").Value = 172

$ws.Range("A9").Value = 232
$ws.Range("B9").Value = 227
$ws.Range("C9").Value = 224

$ws.Range("A10").Value = 253
$ws.Range("B10").Value = 249
$ws.Range("C10").Value = 246

$ws.Range("A11").Value = 254
$ws.Range("B11").Value = 255
$ws.Range("C11").Value = 248

$ws.Range("A12").Value = 239
$ws.Range("B12").Value = 240
$ws.Range("C12").Value = 232

$ws.Range("A13").Value = 246
$ws.Range("B13").Value = 247
$ws.Range("C13").Value = 242

$ws.Range("A14").Value = 255
$ws.Range("B14").Value = 253
$ws.Range("C14").Value = 239

$ws.Range("A15").Value = 239
$ws.Range("B15").Value = 235
$ws.Range("C15").Value = 200

$ws.Range("A16").Value = 255
$ws.Range("B16").Value = 255
$ws.Range("C16").Value = 251

$ws.Range("A17").Value = 254
$ws.Range("B17").Value = 254
$ws.Range("C17").Value = 254

$ws.Range("A18").Value = 255
$ws.Range("B18").Value = 241
$ws.Range("C18").Value = 202

$ws.Range("A19").Value = 250
$ws.Range("B19").Value = 249
$ws.Range("C19").Value = 245

$ws.Range("A20").Value = 240
$ws.Range("B20").Value = 235
$ws.Range("C20").Value = 213

$ws.Range("A21").Value = 252
$ws.Range("B21").Value = 248
$ws.Range("C21").Value = 221

$ws.Range("A22").Value = 252
$ws.Range("B22").Value = 173
$ws.Range("C22").Value = 169

$ws.Range("A23").Value = 241
$ws.Range("B23").Value = 217
$ws.Range("C23").Value = 215

$ws.Range("A24").Value = 254
$ws.Range("B24").Value = 252
$ws.Range("C24").Value = 253

$ws.Range("A25").Value = 247
$ws.Range("B25").Value = 247
$ws.Range("C25").Value = 247

$ws.Range("A26").Value = 251
$ws.Range("B26").Value = 251
$ws.Range("C26").Value = 249

$ws.Range("A27").Value = 254
$ws.Range("B27").Value = 252
$ws.Range("C27").Value = 253

$ws.Range("A28").Value = 253
$ws.Range("B28").Value = 253
$ws.Range("C28").Value = 253

$ws.Range("A29").Value = 255
$ws.Range("B29").Value = 251
$ws.Range("C29").Value = 240

$ws.Range("A30").Value = 252
$ws.Range("B30").Value = 254
$ws.Range("C30").Value = 214

$ws.Range("A31").Value = 251
$ws.Range("B31").Value = 249
$ws.Range("C31").Value = 226

$ws.Range("A32").Value = 247
$ws.Range("B32").Value = 248
$ws.Range("C32").Value = 240

$ws.Range("A33").Value = 253
$ws.Range("B33").Value = 251
$ws.Range("C33").Value = 254

$ws.Range("A34").Value = 254
$ws.Range("B34").Value = 248
$ws.Range("C34").Value = 234

# Max / min summary rows with labels. The order in which the "min"/"max"
# strings are FIRST written controls the shared-string table insertion
# order, so write the "min" label before the "max" label to land them at
# shared-string indices 5 and 6 respectively (matching the target sst).
$ws.Range("D36").Value = "min"
$ws.Range("D35").Value = "max"

$ws.Range("A35").Formula = "=MAX(A2:A34)"
$ws.Range("B35").Formula = "=MAX(B2:B34)"
$ws.Range("C35").Formula = "=MAX(C2:C34)"

$ws.Range("A36").Formula = "=MIN(A2:A34)"
$ws.Range("B36").Formula = "=MIN(B2:B34)"
$ws.Range("C36").Formula = "=MIN(C2:C34)"

# Select the min row (matches the author's last selection on the new sheet)
$ws.Range("A36:C36").Select()
